# Generate Report for Handoff
#
# The localization run moved from "In Translation" to "Ready for handoff"
# and the Xliff/handoff timestamps advanced a little (~40s later). Column
# widths on the Status / generate-date columns are widened so the new,
# longer "Ready for handoff" text is not truncated.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# --- Overview sheet -------------------------------------------------
# E2 = zh-cn status, F2 = de-de status, G2 = Latest HO Xliff Generate Date
$overview.Range("E2").Value = "Ready for handoff"
$overview.Range("F2").Value = "Ready for handoff"
$overview.Range("G2").Value = "2016-08-20 14:43:08"

# widen the (now longer) status columns to fit "Ready for handoff"
$overview.Columns("E").ColumnWidth = 16.3
$overview.Columns("F").ColumnWidth = 16.3

# --- zh-cn sheet ------------------------------------------------------
# C2 = Status, H2 = Latest Handoff Datetime
$zhcn.Range("C2").Value = "Ready for handoff"
$zhcn.Range("H2").Value = "2016-08-20 14:43:03"
$zhcn.Columns("C").ColumnWidth = 16.3

# --- de-de sheet ------------------------------------------------------
# C2 = Status, H2 = Latest Handoff Datetime
$dede.Range("C2").Value = "Ready for handoff"
$dede.Range("H2").Value = "2016-08-20 14:43:08"
$dede.Columns("C").ColumnWidth = 16.3
